# Daily TGP (terminal gate pricing) refresh.
# Rolls the date window forward: each state block keeps its two most-recent
# effective-date rows, so the previous "latest" row becomes the new "previous"
# row, and a brand-new "latest" row of prices is written in its place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New South Wales
$ws.Range("A8").Value = 46025; $ws.Range("D8").Value = 152.21; $ws.Range("E8").Value = 149.96; $ws.Range("F8").Value = 159.96; $ws.Range("G8").Value = 149.97999999999999
$ws.Range("A9").Value = 46025; $ws.Range("D9").Value = 152.21; $ws.Range("E9").Value = 149.96; $ws.Range("F9").Value = 159.96; $ws.Range("G9").Value = 149.97999999999999
$ws.Range("A10").Value = 46025; $ws.Range("D10").Value = 154.46; $ws.Range("E10").Value = 151.83000000000001; $ws.Range("F10").Value = 161.83000000000001; $ws.Range("G10").Value = 152.22999999999999
$ws.Range("A11").Value = 46023; $ws.Range("D11").Value = 152.35; $ws.Range("E11").Value = 150.08000000000001; $ws.Range("F11").Value = 160.08000000000001; $ws.Range("G11").Value = 150.1
$ws.Range("A12").Value = 46023; $ws.Range("D12").Value = 152.35; $ws.Range("E12").Value = 150.08000000000001; $ws.Range("F12").Value = 160.08000000000001; $ws.Range("G12").Value = 150.1
$ws.Range("A13").Value = 46023; $ws.Range("D13").Value = 154.63; $ws.Range("E13").Value = 151.97; $ws.Range("F13").Value = 161.97; $ws.Range("G13").Value = 152.37

# Northern Territory
$ws.Range("A17").Value = 46025; $ws.Range("D17").Value = 158.38; $ws.Range("E17").Value = 155.49; $ws.Range("F17").Value = 165.49
$ws.Range("A18").Value = 46023; $ws.Range("D18").Value = 158.54; $ws.Range("E18").Value = 155.61000000000001; $ws.Range("F18").Value = 165.61

# Queensland
$ws.Range("A22").Value = 46025; $ws.Range("D22").Value = 153.94; $ws.Range("E22").Value = 151.19999999999999; $ws.Range("F22").Value = 160.80000000000001; $ws.Range("G22").Value = 152.28
$ws.Range("A23").Value = 46025; $ws.Range("D23").Value = 158.76; $ws.Range("E23").Value = 156.75; $ws.Range("F23").Value = 166.75
$ws.Range("A24").Value = 46025; $ws.Range("D24").Value = 158.93; $ws.Range("E24").Value = 157.34; $ws.Range("F24").Value = 167.34
$ws.Range("A25").Value = 46025; $ws.Range("D25").Value = 158.91999999999999; $ws.Range("E25").Value = 156.85; $ws.Range("F25").Value = 166.85; $ws.Range("G25").Value = 156.97999999999999
$ws.Range("A26").Value = 46025; $ws.Range("D26").Value = 158.58000000000001; $ws.Range("E26").Value = 158.47; $ws.Range("F26").Value = 168.47
$ws.Range("A27").Value = 46023; $ws.Range("D27").Value = 154.09; $ws.Range("E27").Value = 151.32; $ws.Range("F27").Value = 160.91999999999999; $ws.Range("G27").Value = 152.38999999999999
$ws.Range("A28").Value = 46023; $ws.Range("D28").Value = 158.93; $ws.Range("E28").Value = 156.88; $ws.Range("F28").Value = 166.88
$ws.Range("A29").Value = 46023; $ws.Range("D29").Value = 159.1; $ws.Range("E29").Value = 157.47; $ws.Range("F29").Value = 167.47
$ws.Range("A30").Value = 46023; $ws.Range("D30").Value = 159.09; $ws.Range("E30").Value = 156.99; $ws.Range("F30").Value = 166.99; $ws.Range("G30").Value = 157.11000000000001
$ws.Range("A31").Value = 46023; $ws.Range("D31").Value = 158.74; $ws.Range("E31").Value = 158.61000000000001; $ws.Range("F31").Value = 168.61

# South Australia
$ws.Range("A35").Value = 46025; $ws.Range("D35").Value = 152.02000000000001; $ws.Range("E35").Value = 150.24; $ws.Range("F35").Value = 159.24
$ws.Range("A36").Value = 46023; $ws.Range("D36").Value = 152.19; $ws.Range("E36").Value = 150.38; $ws.Range("F36").Value = 159.38

# Tasmania
$ws.Range("A40").Value = 46025; $ws.Range("D40").Value = 159.43; $ws.Range("E40").Value = 157.57; $ws.Range("F40").Value = 167.57
$ws.Range("A41").Value = 46025; $ws.Range("D41").Value = 159.15; $ws.Range("E41").Value = 157.99; $ws.Range("F41").Value = 167.99
$ws.Range("A42").Value = 46023; $ws.Range("D42").Value = 159.61000000000001; $ws.Range("E42").Value = 157.91999999999999; $ws.Range("F42").Value = 167.92
$ws.Range("A43").Value = 46023; $ws.Range("D43").Value = 159.33000000000001; $ws.Range("E43").Value = 158.34; $ws.Range("F43").Value = 168.34

# Victoria
$ws.Range("A47").Value = 46025; $ws.Range("D47").Value = 153; $ws.Range("E47").Value = 151.9; $ws.Range("F47").Value = 161.9
$ws.Range("A48").Value = 46025; $ws.Range("D48").Value = 152.66999999999999; $ws.Range("E48").Value = 151.86000000000001; $ws.Range("F48").Value = 161.86000000000001
$ws.Range("A49").Value = 46023; $ws.Range("D49").Value = 153.07; $ws.Range("E49").Value = 152.43; $ws.Range("F49").Value = 162.43
$ws.Range("A50").Value = 46023; $ws.Range("D50").Value = 152.74; $ws.Range("E50").Value = 152.38; $ws.Range("F50").Value = 162.38

# Western Australia
$ws.Range("A54").Value = 46025; $ws.Range("D54").Value = 168.77; $ws.Range("E54").Value = 165.32; $ws.Range("F54").Value = 175.32
$ws.Range("A55").Value = 46025; $ws.Range("D55").Value = 161.76; $ws.Range("E55").Value = 163.56; $ws.Range("F55").Value = 173.56
$ws.Range("A56").Value = 46025; $ws.Range("D56").Value = 158.56
$ws.Range("A57").Value = 46025; $ws.Range("D57").Value = 159.34; $ws.Range("E57").Value = 157.97999999999999
$ws.Range("A58").Value = 46025; $ws.Range("D58").Value = 155.11000000000001; $ws.Range("E58").Value = 153.88; $ws.Range("F58").Value = 163.88
$ws.Range("A59").Value = 46025; $ws.Range("D59").Value = 161.30000000000001; $ws.Range("E59").Value = 163.89
$ws.Range("A60").Value = 46023; $ws.Range("D60").Value = 168.93; $ws.Range("E60").Value = 165.5; $ws.Range("F60").Value = 175.5
$ws.Range("A61").Value = 46023; $ws.Range("D61").Value = 161.94; $ws.Range("E61").Value = 163.68; $ws.Range("F61").Value = 173.68
$ws.Range("A62").Value = 46023; $ws.Range("D62").Value = 158.72999999999999
$ws.Range("A63").Value = 46023; $ws.Range("D63").Value = 159.5; $ws.Range("E63").Value = 158.11000000000001
$ws.Range("A64").Value = 46023; $ws.Range("D64").Value = 155.27000000000001; $ws.Range("E64").Value = 154; $ws.Range("F64").Value = 164
$ws.Range("A65").Value = 46023; $ws.Range("D65").Value = 161.46; $ws.Range("E65").Value = 164.05
